# Fruta / hortaliza, semanal
# Insert two new weekly price rows at the top of the data table (row 179),
# pushing the existing rows 179-212 down to 181-214, then populate the two
# new rows with the new week's figures (same Market/Product/Quality/Volume
# as the rows they precede, new date + price columns).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above the current row 179 (this shifts old rows
# 179..212 down to 181..214, matching Excel's native Rows.Insert behaviour).
$ws.Rows.Item(179).Insert()
$ws.Rows.Item(179).Insert()

# --- New row 179 ("Primera" quality) ---
$ws.Range("A179").Value = 4
$ws.Range("B179").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C179").Value = "Los Lagos"
$ws.Range("D179").Value = 44641
$ws.Range("E179").Value = 10
$ws.Range("F179").Value = "Fruta"
$ws.Range("G179").Value = 100101
$ws.Range("H179").Value = "Berries"
$ws.Range("I179").Value = 100112025
$ws.Range("J179").Value = "Frutilla"
$ws.Range("K179").Value = "Sin especificar"
$ws.Range("L179").Value = "Primera"
$ws.Range("M179").Value = 200
$ws.Range("N179").Value = 8000
$ws.Range("O179").Value = 8500
$ws.Range("P179").Value = 8250
$ws.Range("Q179").Value = "$/caja 7 kilos"
$ws.Range("R179").Value = "Región de La Araucanía"
$ws.Range("S179").Value = 1179
$ws.Range("T179").Value = 7

# --- New row 180 ("Segunda" quality) ---
$ws.Range("A180").Value = 4
$ws.Range("B180").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C180").Value = "Los Lagos"
$ws.Range("D180").Value = 44641
$ws.Range("E180").Value = 10
$ws.Range("F180").Value = "Fruta"
$ws.Range("G180").Value = 100101
$ws.Range("H180").Value = "Berries"
$ws.Range("I180").Value = 100112025
$ws.Range("J180").Value = "Frutilla"
$ws.Range("K180").Value = "Sin especificar"
$ws.Range("L180").Value = "Segunda"
$ws.Range("M180").Value = 100
$ws.Range("N180").Value = 5000
$ws.Range("O180").Value = 5000
$ws.Range("P180").Value = 5000
$ws.Range("Q180").Value = "$/caja 7 kilos"
$ws.Range("R180").Value = "Región de La Araucanía"
$ws.Range("S180").Value = 714
$ws.Range("T180").Value = 7

# Apply the same date-time number format used by the rest of column D to
# keep the inserted cells' appearance consistent with the table.
$ws.Range("D179:D180").NumberFormat = $ws.Range("D181").NumberFormat
